# Update countries & provincias Spain
# Refresh the COVID-19 figures table ("Pais" sheet) with the latest data
# snapshot, update the "last updated" timestamp, and swap the ranking of
# "Costa de Marfil" / "Consejo Danes para los Refugiados" (the latter
# overtakes the former in total cases).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 09:25"

# --- Swap the two countries that changed rank order ---
# Row 83 used to be "Costa de Marfil" and row 84 "Consejo Danes para los
# Refugiados". With the new data, Consejo Danes now has more cases than
# Costa de Marfil, so they swap places; Costa de Marfil keeps its old
# (unchanged) figures, now one row further down.
$ws.Range("A83").Value = "Consejo Danes para los Refugiados"
$ws.Range("A84").Value = "Costa de Marfil"

# --- Refresh numeric figures (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected
#     rows ---

# Row 13: India
$ws.Range("B13").Value = 158613
$ws.Range("C13").Value = 527
$ws.Range("D13").Value = 67753
$ws.Range("E13").Value = 86320
$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 4540

# Row 29: Singapur
$ws.Range("B29").Value = 33249
$ws.Range("C29").Value = 373
$ws.Range("E29").Value = 15950

# Row 38: Polonia
$ws.Range("D38").Value = 10560
$ws.Range("E38").Value = 10885

# Row 39: Ucrania
$ws.Range("B39").Value = 22382
$ws.Range("C39").Value = 477
$ws.Range("D39").Value = 8439
$ws.Range("E39").Value = 13274
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 669

# Row 41: Rumania
$ws.Range("E41").Value = 5203
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 1229

# Row 42: Israel
$ws.Range("B42").Value = 16809
$ws.Range("C42").Value = 16
$ws.Range("D42").Value = 14602
$ws.Range("E42").Value = 1926

# Row 48: Afganistan
$ws.Range("B48").Value = 13036
$ws.Range("C48").Value = 580
$ws.Range("D48").Value = 1209
$ws.Range("E48").Value = 11592
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 235

# Row 55: Chequia
$ws.Range("B55").Value = 9103
$ws.Range("C55").Value = 17
$ws.Range("D55").Value = 6377
$ws.Range("E55").Value = 2409

# Row 60: Armenia
$ws.Range("B60").Value = 8216
$ws.Range("C60").Value = 442
$ws.Range("D60").Value = 3287
$ws.Range("E60").Value = 4816
$ws.Range("G60").Value = 15
$ws.Range("H60").Value = 113

# Row 83: now "Consejo Danes para los Refugiados" (new figures)
$ws.Range("B83").Value = 2660
$ws.Range("C83").Value = 114
$ws.Range("D83").Value = 381
$ws.Range("E83").Value = 2210
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 69

# Row 84: now "Costa de Marfil" (carries the old Costa de Marfil figures)
$ws.Range("B84").Value = 2556
$ws.Range("D84").Value = 1302
$ws.Range("E84").Value = 1223
$ws.Range("H84").Value = 31

# Row 85: Bulgaria
$ws.Range("B85").Value = 2477
$ws.Range("C85").Value = 17
$ws.Range("D85").Value = 965
$ws.Range("E85").Value = 1378
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 134

# Row 89: El Salvador
$ws.Range("B89").Value = 2194
$ws.Range("C89").Value = 85
$ws.Range("D89").Value = 1002
$ws.Range("E89").Value = 1153

# Row 110: Letonia
$ws.Range("B110").Value = 1061
$ws.Range("C110").Value = 4
$ws.Range("E110").Value = 296
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 24

# Row 126: Georgia
$ws.Range("B126").Value = 738
$ws.Range("C126").Value = 3
$ws.Range("D126").Value = 573
$ws.Range("E126").Value = 153

# Row 128: Jordania
$ws.Range("D128").Value = 486
$ws.Range("E128").Value = 225

# Row 140: Taiwan
$ws.Range("D140").Value = 420
$ws.Range("E140").Value = 14

# Row 141: Estado de Palestina
$ws.Range("B141").Value = 435
$ws.Range("C141").Value = 1
$ws.Range("E141").Value = 67
